$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: nerve-in-cuff note moves from "polyfasc" (Model 1 / column D)
#     to "cuff too small" (Model 3 / column F). Column C becomes a plain
#     styled (green) blank cell like its neighbours.
$ws.Range("C3").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("B15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = "cuff too small"

$ws.Range("D15").Style = "Normal"
$ws.Range("D15").ClearContents()

# --- Row 24: the "polyfasc" note (Model 1 / column D) is removed outright,
#     no replacement elsewhere on the row.
$ws.Range("D24").Style = "Normal"
$ws.Range("D24").ClearContents()

# --- Model Index lookup table: CorTec 200 -> CorTec 300
$ws.Range("B29").Value = "CorTec 300"

# --- Selection / view: drop the scrolled topLeftCell and move the
#     active selection to C25.
$ws.Range("C25").Select()
